$d = $word.ActiveDocument
$sec = $d.Sections.Item(1)

# wdHeaderFooterPrimary = 1, wdHeaderFooterFirstPage = 2

# Footers: the "PearsonLogo" image (image2.png -> image1.png)
# Footers.Item(1) (Primary)   -> word/footer2.xml (docPr id="4")
# Footers.Item(2) (FirstPage) -> word/footer1.xml (docPr id="2")
$fPrimary = $sec.Footers.Item(1)
if ($fPrimary.Exists -and $fPrimary.Range.InlineShapes.Count -gt 0) {
    $fPrimary.Range.InlineShapes.Item(1).Name = "image1.png"
}
$fFirst = $sec.Footers.Item(2)
if ($fFirst.Exists -and $fFirst.Range.InlineShapes.Count -gt 0) {
    $fFirst.Range.InlineShapes.Item(1).Name = "image1.png"
}

# Headers: the "BTec_Logo-Orange" image (image1.jpg -> image2.jpg)
# Headers.Item(1) (Primary)   -> word/header2.xml (docPr id="3")
# Headers.Item(2) (FirstPage) -> word/header1.xml (docPr id="1")
$hPrimary = $sec.Headers.Item(1)
if ($hPrimary.Exists -and $hPrimary.Range.InlineShapes.Count -gt 0) {
    $hPrimary.Range.InlineShapes.Item(1).Name = "image2.jpg"
}
$hFirst = $sec.Headers.Item(2)
if ($hFirst.Exists -and $hFirst.Range.InlineShapes.Count -gt 0) {
    $hFirst.Range.InlineShapes.Item(1).Name = "image2.jpg"
}
